$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 4
$ws.Range("F5").Value = 431
$ws.Range("F6").Value = 1261
$ws.Range("F8").Value = 7548
$ws.Range("F11").Value = 2078
$ws.Range("F12").Value = 8197
$ws.Range("F15").Value = 59
$ws.Range("F16").Value = 5607
$ws.Range("F18").Value = 2562
$ws.Range("F19").Value = 1102
$ws.Range("F20").Value = 4581
$ws.Range("F21").Value = 329
$ws.Range("F22").Value = 397
$ws.Range("F25").Value = 475
$ws.Range("F26").Value = 2432
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 2785
$ws.Range("F31").Value = 321
$ws.Range("F33").Value = 271
$ws.Range("F34").Value = 632
$ws.Range("F36").Value = 864
$ws.Range("F37").Value = 1621
$ws.Range("F40").Value = 2605
$ws.Range("F42").Value = 2267
$ws.Range("F44").Value = 24

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 95
$ws.Range("F3").Value = 105
$ws.Range("F4").Value = 37
$ws.Range("F5").Value = 1
$ws.Range("F8").Value = 101

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1306

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1306
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 1261
$ws.Range("F7").Value = 7548
$ws.Range("F10").Value = 2078
$ws.Range("F11").Value = 8197
$ws.Range("F14").Value = 59
$ws.Range("F15").Value = 5607
$ws.Range("F17").Value = 2562
$ws.Range("F18").Value = 1102
$ws.Range("F19").Value = 4581
$ws.Range("F20").Value = 397
$ws.Range("F22").Value = 95
$ws.Range("F24").Value = 105
$ws.Range("F25").Value = 475
$ws.Range("F26").Value = 2432
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 2785
$ws.Range("F30").Value = 321
$ws.Range("F32").Value = 271
$ws.Range("F33").Value = 37
$ws.Range("F34").Value = 632
$ws.Range("F35").Value = 1
$ws.Range("F37").Value = 864
$ws.Range("F39").Value = 1621
$ws.Range("F42").Value = 2605
$ws.Range("F45").Value = 2267
$ws.Range("F47").Value = 24
$ws.Range("F49").Value = 101
